# Updated bus voltage magnitude results for Case_2_186 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @(
    @(2, 1.02, 1.028443085198718, 1.033123255901119, 1.051199846960221, 1.056321057559269, 1.03428480327816, 1.033595376084792, 1.035926511816502, 1.053952009755019, 1.059059089846366, 1.01515093291506),
    @(3, 1.02, 1.029658276158951, 1.034047354632595, 1.052415821133193, 1.057599095744064, 1.034584427715368, 1.034449816657127, 1.036659217644161, 1.05497954492784, 1.060149561717304, 1.015438401873628),
    @(4, 1.02, 1.030443747351015, 1.034644271461641, 1.053202924754143, 1.058426086631116, 1.034776198656056, 1.035001383112741, 1.037131669499758, 1.055644149565556, 1.060854630338649, 1.015623848148415),
    @(5, 1.02, 1.030773761398049, 1.034894968269392, 1.053533893799345, 1.058773759324046, 1.034856315475438, 1.035232948741405, 1.037329892568565, 1.055923484253699, 1.061150913919091, 1.015701674884904),
    @(6, 1.02, 1.030829160674332, 1.034937046915187, 1.053589469160404, 1.058832135498467, 1.034869737920652, 1.035271811326516, 1.03736315192488, 1.05597038201074, 1.06120065382834, 1.01571473444078),
    @(7, 1.02, 1.030448157789658, 1.034647622253848, 1.053207346898535, 1.058430732225432, 1.034777271159426, 1.035004478531558, 1.03713431971945, 1.055647882302935, 1.060858589792059, 1.015624888602473),
    @(8, 1.02, 1.028853939958326, 1.033435774803507, 1.051610732882824, 1.056752974826074, 1.034386499159746, 1.033884411260262, 1.036174477163669, 1.054299328981728, 1.059427732801794, 1.015248201743099),
    @(9, 1.02, 1.026038178879723, 1.03129234912249, 1.048799404565999, 1.053796565556439, 1.033681754625613, 1.031900575496479, 1.03447036072836, 1.051920786675179, 1.056902161016236, 1.01458007993028),
    @(10, 1.02, 1.024156429456433, 1.029857942897388, 1.046926477602459, 1.051825493826484, 1.033201027889643, 1.030571100465203, 1.033325630814996, 1.050333486082064, 1.055215486583788, 1.014131708734251),
    @(11, 1.02, 1.023340486264152, 1.029235516257294, 1.046115757075457, 1.050971934326917, 1.032990274653596, 1.029993758205745, 1.032827879372177, 1.049645760526296, 1.05448440761185, 1.013936850558652),
    @(12, 1.02, 1.023037234840678, 1.029004119332248, 1.045814657240171, 1.050654870555295, 1.032911600757733, 1.029779054382744, 1.032642678738131, 1.049390244116625, 1.054212738944226, 1.013864364149877),
    @(13, 1.02, 1.023102291221789, 1.029053763834977, 1.045879242504992, 1.050722882480184, 1.032928494271667, 1.029825120574756, 1.032682419120953, 1.049445056210942, 1.054271017967503, 1.01387991760169),
    @(14, 1.02, 1.023315422972686, 1.029216393004551, 1.046090867311656, 1.050945726034274, 1.032983779421066, 1.029976015898686, 1.032812577037031, 1.04962464078981, 1.054461953722811, 1.013930861004402),
    @(15, 1.02, 1.023446717343561, 1.029316567654368, 1.046221261279759, 1.051083025391619, 1.033017790621639, 1.030068953870649, 1.032892730020115, 1.049735280156848, 1.054579580442398, 1.013962234663846),
    @(16, 1.02, 1.024210556741218, 1.029899223358912, 1.046980287905772, 1.051882139987224, 1.033214960135174, 1.030609381381792, 1.033358621057499, 1.050379119242706, 1.055263990126224, 1.01414462580393),
    @(17, 1.02, 1.024689387118955, 1.030264353450179, 1.047456475511544, 1.052383382034635, 1.033337943918385, 1.030947928425339, 1.033650305230768, 1.050782870030055, 1.055693102917718, 1.014258844262523),
    @(18, 1.02, 1.024968571673268, 1.030477200484097, 1.047734253836157, 1.052675741248743, 1.033409427918643, 1.031145236129715, 1.033820239503816, 1.051018331254639, 1.055943326134091, 1.014325397450581),
    @(19, 1.02, 1.025063747991634, 1.03054975430403, 1.047828973579421, 1.052775427120115, 1.033433759674842, 1.031212485705114, 1.03387814879829, 1.051098610770104, 1.056028633815296, 1.014348078767784),
    @(20, 1.02, 1.024638024462999, 1.030225191614438, 1.047405382402965, 1.052329604231286, 1.033324774824747, 1.030911622207324, 1.033619030997268, 1.050739555524019, 1.055647070541156, 1.014246596793938),
    @(21, 1.02, 1.023252665810183, 1.029168508301497, 1.046028548056199, 1.050880104565167, 1.032967510116149, 1.029931587993981, 1.032774257445211, 1.049571759407121, 1.054405731052291, 1.013915862412328),
    @(22, 1.02, 1.02238062857041, 1.028502971790281, 1.045163095449846, 1.049968665492684, 1.032740622375945, 1.029313935546284, 1.032241299611404, 1.048837145200207, 1.05362459436266, 1.013707294827772),
    @(23, 1.02, 1.022843008671583, 1.028855895619685, 1.045621868470626, 1.050451845244382, 1.032861114479607, 1.029641504491662, 1.032524003276994, 1.049226614363119, 1.054038752854793, 1.013817919608132),
    @(24, 1.02, 1.024661233373348, 1.030242887555278, 1.047428469098648, 1.05235390412674, 1.033330726144885, 1.030928027922074, 1.033633163095695, 1.050759127607739, 1.055667870810085, 1.014252131109481),
    @(25, 1.02, 1.026766914193903, 1.031847431724488, 1.049525963310387, 1.054560881110589, 1.033865865996112, 1.032414656212127, 1.034912434838374, 1.052535971522301, 1.05755559447903, 1.014753324274017)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}
